$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.341.79'
$ws.Range('E2').Value = '  -2.71%  '
$ws.Range('D3').Value = '2.422.68'
$ws.Range('E3').Value = '  -3.39%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '512.71'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -3.82%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '129.06'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -4.58%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.548'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -3.12%  '
$ws.Range('D9').Value = '2.430.32'
$ws.Range('E9').Value = '  -3.27%  '
$ws.Range('B10').Value = 'TRON'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.156'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.27%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0956'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -5.61%  '
$ws.Range('E12').Value = '  -4.26%  '
$ws.Range('E13').Value = '  -3.68%  '
$ws.Range('D14').Value = '2.851.36'
$ws.Range('E14').Value = '  -3.50%  '
$ws.Range('D15').Value = '57.271.57'
$ws.Range('E15').Value = '  -2.66%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '21.48'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -5.79%  '
$ws.Range('E17').Value = '  -4.66%  '
$ws.Range('D18').Value = '2.423.48'
$ws.Range('E18').Value = '  -3.71%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.37'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -5.92%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '314.31'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.68%  '
$ws.Range('E21').Value = '  -4.07%  '
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.64'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -5.05%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '63.44'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -2.40%  '
$ws.Range('E25').Value = '  -4.20%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('E27').Value = '  -2.83%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.21'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -4.60%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '169.17'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.20%  '
$ws.Range('D30').Value = '0.0₃0718'
$ws.Range('E30').Value = '  -5.68%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.66'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -5.03%  '
$ws.Range('B32').Value = 'Aptos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.18'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -4.64%  '
$ws.Range('E33').Value = '  +1.43%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.998'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.998'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('E36').Value = '  -4.01%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.28'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -6.24%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.86'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -4.34%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '36.27'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -2.12%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.45'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -4.82%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.765'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -3.99%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.37'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -5.52%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '264.94'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -5.80%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '4.87'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.59%  '
$ws.Range('E45').Value = '  -3.73%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '120.79'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -6.93%  '
$ws.Range('E47').Value = '  -2.60%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0481'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -3.78%  '
$ws.Range('E49').Value = '  -3.71%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '16.48'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -4.46%  '
$ws.Range('D51').Value = '1.690.65'
$ws.Range('E51').Value = '  -3.70%  '
